$d = $word.ActiveDocument

# The site rebuild dropped the per-page boilerplate that used to follow the
# "Requisitos" section: a blank separator paragraph, the "Ver no Jupiter
# Salvar em pdf Salvar em docx" action line, and the site footer/copyright
# line. Locate the paragraph holding the last requirement line and remove
# the three paragraphs that immediately follow it.
$target = "LOM3223: Materiais e Dispositivos Magnéticos e Supercondutores (Requisito)"

$count = $d.Paragraphs.Count
$startPara = 0
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text -replace "[\r\a]+$", ""
    if ($text -eq $target) {
        $startPara = $i
        break
    }
}

if ($startPara -eq 0) {
    throw "Could not find target paragraph containing: $target"
}

$blank = $d.Paragraphs.Item($startPara + 1)
$jupiter = $d.Paragraphs.Item($startPara + 2)
$footer = $d.Paragraphs.Item($startPara + 3)

# Sanity-check the paragraphs about to be removed so we fail loudly instead
# of silently deleting the wrong content if the document shape differs from
# what's expected.
$blankText = $blank.Range.Text -replace "[\r\a]+$", ""
$jupiterText = $jupiter.Range.Text -replace "[\r\a]+$", ""
$footerText = $footer.Range.Text -replace "[\r\a]+$", ""

if ($blankText -ne "") {
    throw "Unexpected content where blank paragraph was expected: $blankText"
}
if ($jupiterText -notlike "Ver no Jupiter*") {
    throw "Unexpected content where 'Ver no Jupiter...' paragraph was expected: $jupiterText"
}
if ($footerText -notlike "*Contact: luizeleno@usp.br*") {
    throw "Unexpected content where footer paragraph was expected: $footerText"
}

$r = $d.Range($blank.Range.Start, $footer.Range.End)
$r.Delete()
